# Generate Report for Handoff
# Updates the localization-status report: status moves from "Handed back: in
# sync with en-US" to "Ready for handoff", and the generation timestamps are
# refreshed. Column widths on the affected "Status" columns are re-fit to the
# new (shorter) text.

$wb = $excel.ActiveWorkbook

# Old status text being replaced everywhere: "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# Target "Status" column width (29.9777... -> 17.2159881591797 in the saved
# OOXML) is narrower now that "Ready for handoff" is shorter than the old
# status text. The ColumnWidth setter here is quantized to 1/6-character
# steps, so 16.3 is the closest input that lands on the nearest achievable
# stored width (17.166666666666668).
$statusColWidth = 16.3

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-24 22:59:55"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-24 22:59:50"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-24 22:59:55"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
